# Update "want to go" counts (column F) on the individual-category sheets
# as well as the aggregated "全部类型" sheet, matching the upstream data
# refresh recorded in the commit "Update gh-pages to output generated at
# 456a3b4".

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 296
$wsExpo.Range("F14").Value = 7153

# 演出 (Performances) sheet
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 13

# 全部类型 (All categories) sheet - aggregated view of all rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 296
$wsAll.Range("F6").Value = 13
$wsAll.Range("F17").Value = 7153
